$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 9 and 10: status changes from "in Arbeit" to "done" ---
# Copy the "done" cell format (style index used by e.g. B2) onto B9/B10 and set the value.
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B9").Value2 = "done"

$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B10").Value2 = "done"

# --- Add new row 15 for the new TODO item ---
$ws.Range("A15").Value2 = "Generierung eines Order-ID-Hashes, der als eindeutige Order-ID genutzt wird -> darauf aufbauend Verhindern der Trennung von Orders, weil bei jedem Produkt Bestand geprüft wird"

$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B15").Value2 = "done"

$ws.Range("C15").Value2 = "Jonas"

$excel.CutCopyMode = 0

# --- Update the active selection to match the author's last position ---
$ws.Range("B15").Select()
